$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace/extend the data block to 6 rows x 30 columns (A1:AD6).
# Built as a true .NET 2D array so the single Range.Value assignment lands
# on the correct row/column grid (jagged PS arrays are not accepted here).
$arr = New-Object "object[,]" 6,30

$arr[0,0] = 15.78
$arr[0,1] = 22.91
$arr[0,2] = 105.7
$arr[0,3] = 782.6
$arr[0,4] = 0.11550000000000001
$arr[0,5] = 0.17519999999999999
$arr[0,6] = 0.21329999999999999
$arr[0,7] = [double]"9.4789999999999999E-2"
$arr[0,8] = 0.20960000000000001
$arr[0,9] = [double]"7.331E-2"
$arr[0,10] = 0.55200000000000005
$arr[0,11] = 1.0720000000000001
$arr[0,12] = 3.5979999999999999
$arr[0,13] = 58.63
$arr[0,14] = [double]"8.6990000000000001E-3"
$arr[0,15] = [double]"3.9759999999999997E-2"
$arr[0,16] = [double]"5.9499999999999997E-2"
$arr[0,17] = [double]"1.3899999999999999E-2"
$arr[0,18] = [double]"1.495E-2"
$arr[0,19] = [double]"5.9839999999999997E-3"
$arr[0,20] = 20.190000000000001
$arr[0,21] = 30.5
$arr[0,22] = 130.30000000000001
$arr[0,23] = 1272
$arr[0,24] = 0.1855
$arr[0,25] = 0.49249999999999999
$arr[0,26] = 0.73560000000000003
$arr[0,27] = 0.2034
$arr[0,28] = 0.32740000000000002
$arr[0,29] = 0.12520000000000001
$arr[1,0] = 17.95
$arr[1,1] = 20.010000000000002
$arr[1,2] = 114.2
$arr[1,3] = 982
$arr[1,4] = [double]"8.4019999999999997E-2"
$arr[1,5] = [double]"6.7220000000000002E-2"
$arr[1,6] = [double]"7.2929999999999995E-2"
$arr[1,7] = [double]"5.5960000000000003E-2"
$arr[1,8] = 0.21290000000000001
$arr[1,9] = [double]"5.0250000000000003E-2"
$arr[1,10] = 0.55059999999999998
$arr[1,11] = 1.214
$arr[1,12] = 3.3570000000000002
$arr[1,13] = 54.04
$arr[1,14] = [double]"4.0239999999999998E-3"
$arr[1,15] = [double]"8.4220000000000007E-3"
$arr[1,16] = [double]"2.291E-2"
$arr[1,17] = [double]"9.8630000000000002E-3"
$arr[1,18] = [double]"5.0139999999999997E-2"
$arr[1,19] = [double]"1.902E-3"
$arr[1,20] = 20.58
$arr[1,21] = 27.83
$arr[1,22] = 129.19999999999999
$arr[1,23] = 1261
$arr[1,24] = 0.1072
$arr[1,25] = 0.1202
$arr[1,26] = 0.22489999999999999
$arr[1,27] = 0.11849999999999999
$arr[1,28] = 0.48820000000000002
$arr[1,29] = [double]"6.1109999999999998E-2"
$arr[2,0] = 11.41
$arr[2,1] = 10.82
$arr[2,2] = 73.34
$arr[2,3] = 403.3
$arr[2,4] = [double]"9.3729999999999994E-2"
$arr[2,5] = [double]"6.6850000000000007E-2"
$arr[2,6] = [double]"3.5119999999999998E-2"
$arr[2,7] = [double]"2.623E-2"
$arr[2,8] = 0.16669999999999999
$arr[2,9] = [double]"6.1129999999999997E-2"
$arr[2,10] = 0.14080000000000001
$arr[2,11] = 0.4607
$arr[2,12] = 1.103
$arr[2,13] = 10.5
$arr[2,14] = [double]"6.0400000000000002E-3"
$arr[2,15] = [double]"1.529E-2"
$arr[2,16] = [double]"1.5140000000000001E-2"
$arr[2,17] = [double]"6.4599999999999996E-3"
$arr[2,18] = [double]"1.3440000000000001E-2"
$arr[2,19] = [double]"2.2060000000000001E-3"
$arr[2,20] = 12.82
$arr[2,21] = 15.97
$arr[2,22] = 83.74
$arr[2,23] = 510.5
$arr[2,24] = 0.15479999999999999
$arr[2,25] = 0.23899999999999999
$arr[2,26] = 0.2102
$arr[2,27] = [double]"8.9580000000000007E-2"
$arr[2,28] = 0.30159999999999998
$arr[2,29] = [double]"8.523E-2"
$arr[3,0] = 18.66
$arr[3,1] = 17.12
$arr[3,2] = 121.4
$arr[3,3] = 1077
$arr[3,4] = 0.10539999999999999
$arr[3,5] = 0.11
$arr[3,6] = 0.1457
$arr[3,7] = [double]"8.6650000000000005E-2"
$arr[3,8] = 0.1966
$arr[3,9] = [double]"6.2129999999999998E-2"
$arr[3,10] = 0.71279999999999999
$arr[3,11] = 1.581
$arr[3,12] = 4.8949999999999996
$arr[3,13] = 90.47
$arr[3,14] = [double]"8.1019999999999998E-3"
$arr[3,15] = [double]"2.1010000000000001E-2"
$arr[3,16] = [double]"3.3419999999999998E-2"
$arr[3,17] = [double]"1.601E-2"
$arr[3,18] = [double]"2.0449999999999999E-2"
$arr[3,19] = [double]"4.5700000000000003E-3"
$arr[3,20] = 22.25
$arr[3,21] = 24.9
$arr[3,22] = 145.4
$arr[3,23] = 1549
$arr[3,24] = 0.15029999999999999
$arr[3,25] = 0.2291
$arr[3,26] = 0.32719999999999999
$arr[3,27] = 0.16739999999999999
$arr[3,28] = 0.28939999999999999
$arr[3,29] = [double]"8.4559999999999996E-2"
$arr[4,0] = 24.25
$arr[4,1] = 20.2
$arr[4,2] = 166.2
$arr[4,3] = 1761
$arr[4,4] = 0.1447
$arr[4,5] = 0.28670000000000001
$arr[4,6] = 0.42680000000000001
$arr[4,7] = 0.20119999999999999
$arr[4,8] = 0.26550000000000001
$arr[4,9] = [double]"6.8769999999999998E-2"
$arr[4,10] = 1.5089999999999999
$arr[4,11] = 3.12
$arr[4,12] = 9.8070000000000004
$arr[4,13] = 233
$arr[4,14] = [double]"2.333E-2"
$arr[4,15] = [double]"9.8059999999999994E-2"
$arr[4,16] = 0.1278
$arr[4,17] = [double]"1.822E-2"
$arr[4,18] = [double]"4.5469999999999997E-2"
$arr[4,19] = [double]"9.8750000000000001E-3"
$arr[4,20] = 26.02
$arr[4,21] = 23.99
$arr[4,22] = 180.9
$arr[4,23] = 2073
$arr[4,24] = 0.1696
$arr[4,25] = 0.4244
$arr[4,26] = 0.58030000000000004
$arr[4,27] = 0.2248
$arr[4,28] = 0.32219999999999999
$arr[4,29] = [double]"8.0089999999999995E-2"
$arr[5,0] = 14.5
$arr[5,1] = 10.89
$arr[5,2] = 94.28
$arr[5,3] = 640.70000000000005
$arr[5,4] = 0.1101
$arr[5,5] = 0.1099
$arr[5,6] = [double]"8.8419999999999999E-2"
$arr[5,7] = [double]"5.7779999999999998E-2"
$arr[5,8] = 0.18559999999999999
$arr[5,9] = [double]"6.4019999999999994E-2"
$arr[5,10] = 0.29289999999999999
$arr[5,11] = 0.85699999999999998
$arr[5,12] = 1.9279999999999999
$arr[5,13] = 24.19
$arr[5,14] = [double]"3.8180000000000002E-3"
$arr[5,15] = [double]"1.2760000000000001E-2"
$arr[5,16] = [double]"2.8819999999999998E-2"
$arr[5,17] = [double]"1.2E-2"
$arr[5,18] = [double]"1.9099999999999999E-2"
$arr[5,19] = [double]"2.8080000000000002E-3"
$arr[5,20] = 15.7
$arr[5,21] = 15.98
$arr[5,22] = 102.8
$arr[5,23] = 745.5
$arr[5,24] = 0.1313
$arr[5,25] = 0.17879999999999999
$arr[5,26] = 0.25600000000000001
$arr[5,27] = 0.1221
$arr[5,28] = 0.28889999999999999
$arr[5,29] = [double]"8.0060000000000006E-2"

$ws.Range("A1:AD6").Value = $arr

# Match the saved view state: active selection moves to D9.
[void]$ws.Range("D9").Select()
